# Server.xlsx: "unify the conception of DataNode, DataTable, Entity."
#
# The sheet formerly called "Property1" becomes "DataNode" - that rename is
# the actual semantic edit behind this commit. The rest of the diff is the
# usual metadata churn produced when a workbook last saved by Mac Excel gets
# re-opened/re-saved by a newer Windows Excel build: new xr/xr2/xr3
# revision-tracking namespaces, slightly re-measured row heights/column
# widths, phonetic-guide font bookkeeping, a relocated cursor/selection,
# the "Normal" cell style getting its zh-CN display name, etc. We reproduce
# whatever subset of that is actually reachable through the Excel object
# model; the rest is version/application chrome that isn't backed by a
# settable COM property.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- the content edit: rename the worksheet ---
$ws.Name = "DataNode"

# --- incidental re-layout that came along with the re-save ---

# cursor ended up on C36 instead of A9 by the time the file was saved
[void]$ws.Range("C36").Select()

# header rows reflowed slightly (ht 28 -> 27, ht 42 -> 40.5)
$ws.Rows.Item(1).RowHeight = 27
$ws.Rows.Item(8).RowHeight = 40.5

# column widths nudged by the new renderer's font metrics
$ws.Columns.Item(1).ColumnWidth = 27.428571428571427
$ws.Columns.Item(3).ColumnWidth = 14.857142857142858
$ws.Columns.Item(4).ColumnWidth = 10.857142857142858
$ws.Columns.Item(5).ColumnWidth = 7.428571428571429
$ws.Columns.Item(7).ColumnWidth = 21.857142857142858
$ws.Columns.Item(8).ColumnWidth = 8.714285714285714

# "Normal" cell style is relabelled "常规" by zh-CN Excel
$normal = $wb.Styles.Item("Normal")
$normal.Name = "常规"
